$wb = $excel.ActiveWorkbook

# Update "想去人数" (interested-people count) figures on both the
# "展览" sheet and the aggregated "全部类型" sheet, which carry the
# same rows and therefore need to stay in sync.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F5").Value = 4779
    $ws.Range("F9").Value = 737
}
